$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing rows 23-29 (shrinks used range from 29 to 22 rows)
$ws.Range("A23:E29").EntireRow.Delete()

# Update header row (C1, D1, E1 get reordered/relabeled string references)
$ws.Cells.Item(1,3).Value = "最近一次充电结束时间"
$ws.Cells.Item(1,4).Value = "截止一直未充电时间"
$ws.Cells.Item(1,5).Value = "截止一直未充电时长(小时)"

# Update data rows 2-22 with refreshed report data
$ws.Cells.Item(2,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(2,2).Value = "402号直流"
$ws.Cells.Item(2,3).Value = 45915.50368055556
$ws.Cells.Item(2,4).Value = 45930.41322916667
$ws.Cells.Item(2,5).Value = 357.8291666666628

$ws.Cells.Item(3,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(3,2).Value = "101号直流"
$ws.Cells.Item(3,3).Value = 45926.04369212963
$ws.Cells.Item(3,4).Value = 45930.41322916667
$ws.Cells.Item(3,5).Value = 104.86888888897374

$ws.Cells.Item(4,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(4,2).Value = "602号直流"
$ws.Cells.Item(4,3).Value = 45926.24207175926
$ws.Cells.Item(4,4).Value = 45930.41322916667
$ws.Cells.Item(4,5).Value = 100.1077777777682

$ws.Cells.Item(5,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(5,2).Value = "008B号直流"
$ws.Cells.Item(5,3).Value = 45926.52563657407
$ws.Cells.Item(5,4).Value = 45930.41322916667
$ws.Cells.Item(5,5).Value = 93.3022222223226

$ws.Cells.Item(6,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(6,2).Value = "701号直流"
$ws.Cells.Item(6,3).Value = 45927.457337962966
$ws.Cells.Item(6,4).Value = 45930.41322916667
$ws.Cells.Item(6,5).Value = 70.94138888886664

$ws.Cells.Item(7,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(7,2).Value = "702号直流"
$ws.Cells.Item(7,3).Value = 45928.5321875
$ws.Cells.Item(7,4).Value = 45930.41322916667
$ws.Cells.Item(7,5).Value = 45.14500000001863

$ws.Cells.Item(8,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(8,2).Value = "306号直流"
$ws.Cells.Item(8,3).Value = 45928.59892361111
$ws.Cells.Item(8,4).Value = 45930.41322916667
$ws.Cells.Item(8,5).Value = 43.543333333334886

$ws.Cells.Item(9,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(9,2).Value = "903号直流"
$ws.Cells.Item(9,3).Value = 45928.766284722224
$ws.Cells.Item(9,4).Value = 45930.41322916667
$ws.Cells.Item(9,5).Value = 39.5266666666721

$ws.Cells.Item(10,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(10,2).Value = "B05号直流"
$ws.Cells.Item(10,3).Value = 45929.019733796296
$ws.Cells.Item(10,4).Value = 45930.41322916667
$ws.Cells.Item(10,5).Value = 33.44388888892718

$ws.Cells.Item(11,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(11,2).Value = "203号直流"
$ws.Cells.Item(11,3).Value = 45929.02199074074
$ws.Cells.Item(11,4).Value = 45930.41322916667
$ws.Cells.Item(11,5).Value = 33.38972222234588

$ws.Cells.Item(12,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(12,2).Value = "103号直流"
$ws.Cells.Item(12,3).Value = 45929.03383101852
$ws.Cells.Item(12,4).Value = 45930.41322916667
$ws.Cells.Item(12,5).Value = 33.10555555560859

$ws.Cells.Item(13,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(13,2).Value = "801号直流"
$ws.Cells.Item(13,3).Value = 45929.03690972222
$ws.Cells.Item(13,4).Value = 45930.41322916667
$ws.Cells.Item(13,5).Value = 33.031666666734964

$ws.Cells.Item(14,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(14,2).Value = "002B号直流"
$ws.Cells.Item(14,3).Value = 45929.03873842592
$ws.Cells.Item(14,4).Value = 45930.41322916667
$ws.Cells.Item(14,5).Value = 32.98777777788928

$ws.Cells.Item(15,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(15,2).Value = "503号直流"
$ws.Cells.Item(15,3).Value = 45929.158055555556
$ws.Cells.Item(15,4).Value = 45930.41322916667
$ws.Cells.Item(15,5).Value = 30.124166666704696

$ws.Cells.Item(16,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(16,2).Value = "B01号直流"
$ws.Cells.Item(16,3).Value = 45929.52564814815
$ws.Cells.Item(16,4).Value = 45930.41322916667
$ws.Cells.Item(16,5).Value = 21.301944444479886

$ws.Cells.Item(17,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(17,2).Value = "905号直流"
$ws.Cells.Item(17,3).Value = 45929.531018518515
$ws.Cells.Item(17,4).Value = 45930.41322916667
$ws.Cells.Item(17,5).Value = 21.173055555671453

$ws.Cells.Item(18,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(18,2).Value = "402号直流"
$ws.Cells.Item(18,3).Value = 45929.53586805556
$ws.Cells.Item(18,4).Value = 45930.41322916667
$ws.Cells.Item(18,5).Value = 21.05666666664183

$ws.Cells.Item(19,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(19,2).Value = "905号直流"
$ws.Cells.Item(19,3).Value = 45929.56361111111
$ws.Cells.Item(19,4).Value = 45930.41322916667
$ws.Cells.Item(19,5).Value = 20.39083333342569

$ws.Cells.Item(20,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(20,2).Value = "110号直流"
$ws.Cells.Item(20,3).Value = 45929.56958333333
$ws.Cells.Item(20,4).Value = 45930.41322916667
$ws.Cells.Item(20,5).Value = 20.247500000114087

$ws.Cells.Item(21,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(21,2).Value = "006B号直流"
$ws.Cells.Item(21,3).Value = 45929.611122685186
$ws.Cells.Item(21,4).Value = 45930.41322916667
$ws.Cells.Item(21,5).Value = 19.250555555569008

$ws.Cells.Item(22,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(22,2).Value = "B03号直流"
$ws.Cells.Item(22,3).Value = 45929.7096875
$ws.Cells.Item(22,4).Value = 45930.41322916667
$ws.Cells.Item(22,5).Value = 16.88500000006752

# Update selected cell in the sheet view
$null = $ws.Range("G5").Select()
